# Insert a new data row before the existing row 240 (Feria Lagunitas de
# Puerto Montt - Piña weekly price sheet). This pushes the old rows
# 240..326 down to 241..327 and grows the sheet by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 240, shifting everything below it down.
$ws.Rows.Item(240).Insert()

# Populate the newly inserted row 240 with the new weekly record. The
# categorical columns (Mercado/Region/Codreg/Tipo/Producto/Categoria/
# Variedad/Calidad/Unidad/Origen/Kg per unit) repeat the same values as
# the rest of this block, only the date and price/volume figures differ.
$ws.Range("A240").Value = 4
$ws.Range("B240").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C240").Value = "Los Lagos"
$ws.Range("D240").Value = 44876
$ws.Range("E240").Value = 10
$ws.Range("F240").Value = "Fruta"
$ws.Range("G240").Value = 100108
$ws.Range("H240").Value = "Tropicales y subtropicales"
$ws.Range("I240").Value = 100108005
$ws.Range("J240").Value = "Piña"
$ws.Range("K240").Value = "Caramelo"
$ws.Range("L240").Value = "Segunda"
$ws.Range("M240").Value = 200
$ws.Range("N240").Value = 31000
$ws.Range("O240").Value = 32000
$ws.Range("P240").Value = 31500
$ws.Range("Q240").Value = "$/caja 14 unidades"
$ws.Range("R240").Value = "Ecuador"
$ws.Range("S240").Value = 2250
$ws.Range("T240").Value = 14

# Keep the date column's number format consistent with the rest of the
# column (it already inherited style from the row above on insert, but
# set it explicitly to be safe).
$ws.Range("D240").NumberFormat = $ws.Range("D241").NumberFormat
